$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7, pushing existing rows 7-18 down to 8-19.
$ws.Rows.Item(7).Insert()

# Fill in the new row 7 with the new weekly data point.
$ws.Cells.Item(7, 1).Value = 10
$ws.Cells.Item(7, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(7, 3).Value = "La Araucanía"
$ws.Cells.Item(7, 4).Value = 44763
$ws.Cells.Item(7, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(7, 5).Value = 9
$ws.Cells.Item(7, 6).Value = "Fruta"
$ws.Cells.Item(7, 7).Value = 100108
$ws.Cells.Item(7, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(7, 9).Value = 100108001
$ws.Cells.Item(7, 10).Value = "Guayaba"
$ws.Cells.Item(7, 11).Value = "Sin especificar"
$ws.Cells.Item(7, 12).Value = "Primera"
$ws.Cells.Item(7, 13).Value = 50
$ws.Cells.Item(7, 14).Value = 2300
$ws.Cells.Item(7, 15).Value = 2300
$ws.Cells.Item(7, 16).Value = 2300
$ws.Cells.Item(7, 17).Value = "$/kilo"
$ws.Cells.Item(7, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(7, 19).Value = 2300
$ws.Cells.Item(7, 20).Value = 1
